$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect so we can write the updated figures.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A59).
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-10 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) figures for rows 2-56.
$ws.Range("D2").Value = 0.01978569641532148
$ws.Range("E2").Value = -0.003650128430444766
$ws.Range("D3").Value = 0.01981511951395181
$ws.Range("E3").Value = -0.02717782577393812
$ws.Range("D4").Value = 0.01923521698821953
$ws.Range("E4").Value = -0.02514184002669939
$ws.Range("D5").Value = 0.02046563747639688
$ws.Range("E5").Value = -0.0171476369719783
$ws.Range("D6").Value = 0.0195861543100649
$ws.Range("E6").Value = -0.04064241232382837
$ws.Range("D7").Value = 0.007997733173152796
$ws.Range("E7").Value = 0.008026755852842848
$ws.Range("D8").Value = 0.019867189482437
$ws.Range("E8").Value = -0.008464079273327707
$ws.Range("D9").Value = 0.02405739537101547
$ws.Range("E9").Value = 0.0006671114076051143
$ws.Range("D10").Value = 0.02024416178852496
$ws.Range("E10").Value = 0.00179694519317164
$ws.Range("D11").Value = 0.02011684001627008
$ws.Range("E11").Value = -0.02122114668652253
$ws.Range("D12").Value = 0.01354532467854373
$ws.Range("E12").Value = 0.01079515534491837
$ws.Range("D13").Value = 0.01405247190584466
$ws.Range("E13").Value = -0.01690269529465505
$ws.Range("D14").Value = 0.008453167075588299
$ws.Range("E14").Value = 0.00286895620622718
$ws.Range("D15").Value = 0.01436275185503721
$ws.Range("E15").Value = 0.001191895113230146
$ws.Range("D16").Value = 0.020631476759586
$ws.Range("E16").Value = -0.006084806997528003
$ws.Range("D17").Value = 0.02367685662872989
$ws.Range("E17").Value = 0.008857022353437438
$ws.Range("D18").Value = 0.02118926738089307
$ws.Range("E18").Value = 0.002558362647905321
$ws.Range("D19").Value = 0.01967228374423731
$ws.Range("E19").Value = -0.02064014358360755
$ws.Range("D20").Value = 0.01959685361865774
$ws.Range("E20").Value = -0.00847164591977867
$ws.Range("D21").Value = 0.02519455355262518
$ws.Range("E21").Value = 0.003729996390326118
$ws.Range("D22").Value = 0.01669092140484062
$ws.Range("E22").Value = -0.005384615384615321
$ws.Range("D23").Value = 0.02037647657145649
$ws.Range("E23").Value = 0.02058319039451106
$ws.Range("D24").Value = 0.01959132564255144
$ws.Range("E24").Value = -0.02325581395348852
$ws.Range("D25").Value = 0.02048239972652567
$ws.Range("E25").Value = -0.02322787344813759
$ws.Range("D26").Value = 0.0185472514456995
$ws.Range("E26").Value = -0.006489760599942551
$ws.Range("D27").Value = 0.02020528763397095
$ws.Range("E27").Value = -0.02937127122533267
$ws.Range("D28").Value = 0.02126487582828252
$ws.Range("E28").Value = -0.01962264150943405
$ws.Range("D29").Value = 0.02043139968889977
$ws.Range("E29").Value = -0.01881720430107525
$ws.Range("D30").Value = 0.0201714064900936
$ws.Range("E30").Value = -0.01187255785993391
$ws.Range("D31").Value = 0.02119943172405627
$ws.Range("E31").Value = -0.01470353204411057
$ws.Range("D32").Value = 0.02194535185478756
$ws.Range("E32").Value = 0.005590496156533842
$ws.Range("D33").Value = 0.01954692351189112
$ws.Range("E33").Value = 0.00759013282732468
$ws.Range("D34").Value = 0.02006994138027144
$ws.Range("E34").Value = -0.008502963153826482
$ws.Range("D35").Value = 0.0201574973889229
$ws.Range("E35").Value = 0.001167728237791765
$ws.Range("D36").Value = 0.01701190066262602
$ws.Range("E36").Value = 0.001415094339622724
$ws.Range("D37").Value = 0.02059616904122961
$ws.Range("E37").Value = 0.004571428571428449
$ws.Range("D38").Value = 0.01927319953372413
$ws.Range("E38").Value = -0.002498126405196199
$ws.Range("D39").Value = 0.020199224692435
$ws.Range("E39").Value = -0.02447163515016693
$ws.Range("D40").Value = 0.01666845285679564
$ws.Range("E40").Value = -0.006354708261120723
$ws.Range("D41").Value = 0.01331564618741729
$ws.Range("E41").Value = -0.0228197985858154
$ws.Range("D42").Value = 0.01690062785326041
$ws.Range("E42").Value = -0.01671309192200565
$ws.Range("D43").Value = 0.01920436731511016
$ws.Range("E43").Value = -0.02683504340962917
$ws.Range("D44").Value = 0.01303282779694638
$ws.Range("E44").Value = -0.006937033084311595
$ws.Range("D45").Value = 0.01680504736316431
$ws.Range("E45").Value = 0.01315789473684226
$ws.Range("D46").Value = 0.01642700512621707
$ws.Range("E46").Value = -0.001085540599218504
$ws.Range("D47").Value = 0.01377072344623303
$ws.Range("E47").Value = -0.003172588832487166
$ws.Range("D48").Value = 0.02146994590964541
$ws.Range("E48").Value = -0.04709302325581399
$ws.Range("D49").Value = 0.01966515087184208
$ws.Range("E49").Value = -0.01510713735162617
$ws.Range("D50").Value = 0.01808860775068614
$ws.Range("E50").Value = -0.003302509907529783
$ws.Range("D51").Value = 0.01921863305990062
$ws.Range("E51").Value = -0.01519832985386216
$ws.Range("D52").Value = 0.006120004515108228
$ws.Range("E52").Value = 0.02027972027972003
$ws.Range("D53").Value = 0.02169784118267304
$ws.Range("E53").Value = -0.01023192360163716
$ws.Range("D54").Value = 0.01865816761144534
$ws.Range("E54").Value = 0.006097560975609539
$ws.Range("D55").Value = 0.01964981519619233
$ws.Range("E55").Value = -0.03068252974326857
$ws.Range("E56").Value = -0.009463716772182384
